$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns("A")

# Locate the relevant storyText rows by their current contents so the
# script is resilient to exact row-number assumptions.
$kneeRow  = $col.Find("got down on one knee").Row
$meTooRow = $col.Find("Alex, me too").Row
$marryRow = $col.Find("Will you marry me").Row

# Fold the standalone "Will you marry me?" line into the proposal
# sentence, so Alex asks the question as he gets down on one knee.
$ws.Cells.Item($kneeRow, 1).Value = "Seemingly out of nowhere, Alex took Chloe’s hand and he got down on one knee, and asked, “Will you marry me?”"

# Delete the two rows that are no longer needed: the now-redundant
# "Will you marry me?" row and Chloe's "Alex, me too," she replied."
# reply, which no longer fits after the proposal line was rewritten.
# Delete the lower row first so the other row's index doesn't shift.
if ($marryRow -gt $meTooRow) {
    $ws.Rows($marryRow).Delete()
    $ws.Rows($meTooRow).Delete()
} else {
    $ws.Rows($meTooRow).Delete()
    $ws.Rows($marryRow).Delete()
}
